$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 27843
$ws.Range("E2").Value = 543181485383
$ws.Range("F2").Value = 4667714232
$ws.Range("G2").Value = -0.51074

$ws.Range("D3").Value = 1620.97
$ws.Range("E3").Value = 194910373595
$ws.Range("F3").Value = 4555745239
$ws.Range("G3").Value = -1.20362

$ws.Range("D4").Value = 1.001
$ws.Range("E4").Value = 83497079037
$ws.Range("F4").Value = 11283577884
$ws.Range("G4").Value = 0.01407

$ws.Range("D5").Value = 210.27
$ws.Range("E5").Value = 32363234631
$ws.Range("F5").Value = 208839090
$ws.Range("G5").Value = -1.1398

$ws.Range("D6").Value = 0.519411
$ws.Range("E6").Value = 27714558076
$ws.Range("F6").Value = 283173955
$ws.Range("G6").Value = -0.70167

$ws.Range("D7").Value = 0.999915
$ws.Range("E7").Value = 25538586130
$ws.Range("F7").Value = 2248218712
$ws.Range("G7").Value = 0.0305

$ws.Range("D8").Value = 1619.8
$ws.Range("E8").Value = 14240451895
$ws.Range("F8").Value = 42235171
$ws.Range("G8").Value = -1.25206

$ws.Range("D9").Value = 23.22
$ws.Range("E9").Value = 9621133504
$ws.Range("F9").Value = 240271740
$ws.Range("G9").Value = -1.36738

$ws.Range("D10").Value = 0.256061
$ws.Range("E10").Value = 9280973293
$ws.Range("F10").Value = 91551002
$ws.Range("G10").Value = -0.82157

$ws.Range("D11").Value = 0.060984
$ws.Range("E11").Value = 8618302159
$ws.Range("F11").Value = 137713814
$ws.Range("G11").Value = -0.81281

$ws.Range("D12").Value = 0.087643
$ws.Range("E12").Value = 7797204560
$ws.Range("F12").Value = 150078518
$ws.Range("G12").Value = -0.5502

$ws.Range("D13").Value = 2.03
$ws.Range("E13").Value = 7005225406
$ws.Range("F13").Value = 12154480
$ws.Range("G13").Value = -2.93037

$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.5585869999999999
$ws.Range("E14").Value = 5195438559
$ws.Range("F14").Value = 159289213
$ws.Range("G14").Value = -2.86546

$ws.Range("B15").Value = "DOT"
$ws.Range("C15").Value = "Polkadot"
$ws.Range("D15").Value = 4.01
$ws.Range("E15").Value = 5148118672
$ws.Range("F15").Value = 59178587
$ws.Range("G15").Value = -1.9589

$ws.Range("D16").Value = 65.15000000000001
$ws.Range("E16").Value = 4799278724
$ws.Range("F16").Value = 242324233
$ws.Range("G16").Value = -1.06136

$ws.Range("B17").Value = "WBTC"
$ws.Range("C17").Value = "Wrapped Bitcoin"
$ws.Range("D17").Value = 27856
$ws.Range("E17").Value = 4534784820
$ws.Range("F17").Value = 50537413
$ws.Range("G17").Value = -0.53415

$ws.Range("B18").Value = "BCH"
$ws.Range("C18").Value = "Bitcoin Cash"
$ws.Range("D18").Value = 227.86
$ws.Range("E18").Value = 4451834097
$ws.Range("F18").Value = 88738457
$ws.Range("G18").Value = -2.35166

$ws.Range("B19").Value = "SHIB"
$ws.Range("C19").Value = "Shiba Inu"
$ws.Range("D19").Value = 0.00000719
$ws.Range("E19").Value = 4234829252
$ws.Range("F19").Value = 53035715
$ws.Range("G19").Value = -0.6211

$ws.Range("B20").Value = "LINK"
$ws.Range("C20").Value = "Chainlink"
$ws.Range("D20").Value = 7.59
$ws.Range("E20").Value = 4225501918
$ws.Range("F20").Value = 178343111
$ws.Range("G20").Value = -0.30151

$ws.Range("D21").Value = 0.999371
$ws.Range("E21").Value = 3830717813
$ws.Range("F21").Value = 66564302
$ws.Range("G21").Value = -0.08259

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "LEO Token"
$ws.Range("D22").Value = 3.84
$ws.Range("E22").Value = 3572187099
$ws.Range("F22").Value = 94473
$ws.Range("G22").Value = -0.11313

$ws.Range("B23").Value = "AVAX"
$ws.Range("C23").Value = "Avalanche"
$ws.Range("D23").Value = 10.05
$ws.Range("E23").Value = 3560326948
$ws.Range("F23").Value = 175127595
$ws.Range("G23").Value = -4.09623

$ws.Range("B24").Value = "TUSD"
$ws.Range("C24").Value = "TrueUSD"
$ws.Range("D24").Value = 0.999342
$ws.Range("E24").Value = 3414527602
$ws.Range("F24").Value = 92031925
$ws.Range("G24").Value = -0.04549

$ws.Range("B25").Value = "UNI"
$ws.Range("C25").Value = "Uniswap"
$ws.Range("D25").Value = 4.3
$ws.Range("E25").Value = 3243716129
$ws.Range("F25").Value = 35493289
$ws.Range("G25").Value = -1.49434

$ws.Range("D26").Value = 0.110797
$ws.Range("E26").Value = 3075498370
$ws.Range("F26").Value = 32658899
$ws.Range("G26").Value = -0.19885

$ws.Range("D27").Value = 153.94
$ws.Range("E27").Value = 2792048247
$ws.Range("F27").Value = 68984675
$ws.Range("G27").Value = 1.16349

$ws.Range("D28").Value = 42.85
$ws.Range("E28").Value = 2570480471
$ws.Range("F28").Value = 1207154
$ws.Range("G28").Value = -0.19429

$ws.Range("B29").Value = "BUSD"
$ws.Range("C29").Value = "BUSD"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 2214208023
$ws.Range("F29").Value = 1595223858
$ws.Range("G29").Value = -0.00874

$ws.Range("B30").Value = "ETC"
$ws.Range("C30").Value = "Ethereum Classic"
$ws.Range("D30").Value = 15.43
$ws.Range("E30").Value = 2208836903
$ws.Range("F30").Value = 50433624
$ws.Range("G30").Value = -1.49947

$ws.Range("D31").Value = 6.88
$ws.Range("E31").Value = 2008088814
$ws.Range("F31").Value = 122999552
$ws.Range("G31").Value = -1.17377

$ws.Range("D32").Value = 0.04795734
$ws.Range("E32").Value = 1603977505
$ws.Range("F32").Value = 19266014
$ws.Range("G32").Value = -0.80767

$ws.Range("D33").Value = 3.42
$ws.Range("E33").Value = 1555539778
$ws.Range("F33").Value = 67793600
$ws.Range("G33").Value = 1.56894

$ws.Range("D34").Value = 1.57
$ws.Range("E34").Value = 1394637434
$ws.Range("F34").Value = 20067004
$ws.Range("G34").Value = -0.24166

$ws.Range("D35").Value = 3.07
$ws.Range("E35").Value = 1366329278
$ws.Range("F35").Value = 12142992
$ws.Range("G35").Value = -1.70944

$ws.Range("B36").Value = "CRO"
$ws.Range("C36").Value = "Cronos"
$ws.Range("D36").Value = 0.050046
$ws.Range("E36").Value = 1316610464
$ws.Range("F36").Value = 3075579
$ws.Range("G36").Value = -1.01384

$ws.Range("B37").Value = "QNT"
$ws.Range("C37").Value = "Quant"
$ws.Range("D37").Value = 87.52
$ws.Range("E37").Value = 1273274244
$ws.Range("F37").Value = 7477937
$ws.Range("G37").Value = -0.22121

$ws.Range("B38").Value = "MKR"
$ws.Range("C38").Value = "Maker"
$ws.Range("D38").Value = 1383.35
$ws.Range("E38").Value = 1246758994
$ws.Range("F38").Value = 43444629
$ws.Range("G38").Value = -2.18273

$ws.Range("D39").Value = 5.21
$ws.Range("E39").Value = 1246109881
$ws.Range("F39").Value = 26570343
$ws.Range("G39").Value = -1.4916

$ws.Range("B40").Value = "VET"
$ws.Range("C40").Value = "VeChain"
$ws.Range("D40").Value = 0.01694127
$ws.Range("E40").Value = 1230520110
$ws.Range("F40").Value = 22302401
$ws.Range("G40").Value = 0.09877

$ws.Range("B41").Value = "MNT"
$ws.Range("C41").Value = "Mantle"
$ws.Range("D41").Value = 0.382248
$ws.Range("E41").Value = 1221536771
$ws.Range("F41").Value = 21499108
$ws.Range("G41").Value = -2.82987

$ws.Range("B42").Value = "OP"
$ws.Range("C42").Value = "Optimism"
$ws.Range("D42").Value = 1.27
$ws.Range("E42").Value = 1109907699
$ws.Range("F42").Value = 60650998
$ws.Range("G42").Value = -2.46999

$ws.Range("B43").Value = "ARB"
$ws.Range("C43").Value = "Arbitrum"
$ws.Range("D43").Value = 0.843456
$ws.Range("E43").Value = 1075061997
$ws.Range("F43").Value = 68572422
$ws.Range("G43").Value = -4.45274

$ws.Range("B44").Value = "NEAR"
$ws.Range("C44").Value = "NEAR Protocol"
$ws.Range("D44").Value = 1.08
$ws.Range("E44").Value = 1054859374
$ws.Range("F44").Value = 26602837
$ws.Range("G44").Value = -1.27591

$ws.Range("D45").Value = 0.04888521
$ws.Range("E45").Value = 1028373473
$ws.Range("F45").Value = 9156830
$ws.Range("G45").Value = -3.48943

$ws.Range("D46").Value = 65.5
$ws.Range("E46").Value = 955582001
$ws.Range("F46").Value = 55199316
$ws.Range("G46").Value = -1.9766

$ws.Range("D47").Value = 1756.3
$ws.Range("E47").Value = 936515414
$ws.Range("F47").Value = 8380306
$ws.Range("G47").Value = -1.19005

$ws.Range("B48").Value = "ALGO"
$ws.Range("C48").Value = "Algorand"
$ws.Range("D48").Value = 0.100913
$ws.Range("E48").Value = 798367111
$ws.Range("F48").Value = 26889711
$ws.Range("G48").Value = 0.76251

$ws.Range("B49").Value = "GRT"
$ws.Range("C49").Value = "The Graph"
$ws.Range("D49").Value = 0.08459999999999999
$ws.Range("E49").Value = 781364370
$ws.Range("F49").Value = 15997145
$ws.Range("G49").Value = -1.06353

$ws.Range("D50").Value = 5.17
$ws.Range("E50").Value = 756728465
$ws.Range("F50").Value = 3988397
$ws.Range("G50").Value = -0.2926

$ws.Range("D51").Value = 0.999309
$ws.Range("E51").Value = 728779070
$ws.Range("F51").Value = 21387464
$ws.Range("G51").Value = -0.07804999999999999
